$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15, pushing existing rows 15-64 down to 16-65
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row 15 with the new record's data
$ws.Cells.Item(15, 1).Value = 7
$ws.Cells.Item(15, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(15, 3).Value = "Ñuble"
$ws.Cells.Item(15, 4).Value = 44608
$ws.Cells.Item(15, 5).Value = 16
$ws.Cells.Item(15, 6).Value = 100112021
$ws.Cells.Item(15, 7).Value = "Ají"
$ws.Cells.Item(15, 8).Value = "Americana (o)"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 100
$ws.Cells.Item(15, 11).Value = 9000
$ws.Cells.Item(15, 12).Value = 9500
$ws.Cells.Item(15, 13).Value = 9250
$ws.Cells.Item(15, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(15, 15).Value = "Región del Maule"
$ws.Cells.Item(15, 16).Value = 617
$ws.Cells.Item(15, 17).Value = 15
$ws.Cells.Item(15, 18).Value = "Hortaliza"
